$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly price log for "Arándano (blue)" at Vega Monumental
# Concepción, sorted newest-first within each quality grade. A new week's
# observation is inserted at the top of the "Primera" quality block (row 111),
# pushing the existing rows 111-141 down to 112-142.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly data point.
$ws.Cells.Item(111, 1).Value = 11
$ws.Cells.Item(111, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(111, 3).Value = 'Bíobío'
$ws.Cells.Item(111, 4).Value = (Get-Date -Year 2023 -Month 2 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(111, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item(111, 5).Value = 8
$ws.Cells.Item(111, 6).Value = 'Fruta'
$ws.Cells.Item(111, 7).Value = 100101
$ws.Cells.Item(111, 8).Value = 'Berries'
$ws.Cells.Item(111, 9).Value = 100101001
$ws.Cells.Item(111, 10).Value = 'Arándano (blue)'
$ws.Cells.Item(111, 11).Value = 'Sin especificar'
$ws.Cells.Item(111, 12).Value = 'Primera'
$ws.Cells.Item(111, 13).Value = 140
$ws.Cells.Item(111, 14).Value = 3000
$ws.Cells.Item(111, 15).Value = 3500
$ws.Cells.Item(111, 16).Value = 3286
$ws.Cells.Item(111, 17).Value = '$/bandeja 2 kilos'
$ws.Cells.Item(111, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(111, 19).Value = 1643
$ws.Cells.Item(111, 20).Value = 2
